# Apply repulled data values to column F ("dSF") for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    3  = -1
    8  = -3
    11 = 6
    13 = -3
    14 = -4
    15 = 2
    16 = 1
    17 = 1
    18 = 4
    19 = 7
    20 = 6
    21 = 3
    22 = -1
    23 = -1
    24 = -1
    25 = -1
    26 = 5
    27 = 2
    28 = -2
    29 = -6
    30 = -1
    31 = 0
    32 = 3
    33 = -3
    34 = -1
    35 = -5
    38 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
